# This edit reshuffles the per-row data (date, volume, prices, unit,
# origin and $/Kg) among rows 2-15 of the sheet, while A,B,C,E-L,T
# (market/product identifying columns) stay the same for every row.
#
# Effectively it's a permutation of the "data block" (columns D,M,N,O,P,Q,R,S)
# across the existing rows. Mapping: new row -> row whose data block it
# should now contain (derived from matching the unique Fecha/D values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "Q", "R", "S")

# 1) Snapshot all current values for the columns that will be moved around,
#    for every data row (2-15), before any cell gets overwritten.
$orig = @{}
for ($r = 2; $r -le 15; $r++) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $orig[$addr] = $ws.Range($addr).Value2
    }
}

# 2) new row number -> old row number that supplies its data block
$perm = @{
    2  = 12
    3  = 11
    4  = 15
    5  = 9
    6  = 3
    7  = 10
    8  = 5
    9  = 13
    10 = 7
    11 = 2
    12 = 4
    13 = 6
    14 = 14
    15 = 8
}

# 3) Write back the values according to the permutation, using the
#    untouched snapshot captured in step 1.
foreach ($newRow in $perm.Keys) {
    $oldRow = $perm[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $orig["$c$oldRow"]
    }
}
